# Automatische test-sync: 2025-08-08 20:11:50
#
# Appends the newly logged mail (row 5) to the "Logs" sheet, appends the
# matching rollup row (row 4) to the "Dashboard" sheet, and widens the
# conditional-formatting ranges on "Logs" so the new row inherits the same
# rules as the rest of the table.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# Logs!A5:J5 - new mail-log entry
# ---------------------------------------------------------------------
$wsLogs.Cells.Item(5, 1).Value = "Kun jij dit even regelen?"
$wsLogs.Cells.Item(5, 2).Value = "`"Testbedrijf 123 B.V.`" <admin@testbedrijf123.nl>"
$wsLogs.Cells.Item(5, 3).Value = "Testmail #1: Kun jij dit even regelen?`nTestbedrijf 123 B.V."
$wsLogs.Cells.Item(5, 4).Value = "Intern verzoek / Actie voor medewerker"
$wsLogs.Cells.Item(5, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$wsLogs.Cells.Item(5, 6).Value = "2025-08-08 20:11:40"
$wsLogs.Cells.Item(5, 7).Value = "Ja"
$wsLogs.Cells.Item(5, 8).Value = "Ja"
$wsLogs.Cells.Item(5, 9).Value = "Nee"
$wsLogs.Cells.Item(5, 10).Value = "Nee"

# ---------------------------------------------------------------------
# Dashboard!A4:B4 - rollup for the new "Intern verzoek / Actie voor
# medewerker" category
# ---------------------------------------------------------------------
$wsDash.Cells.Item(4, 1).Value = "Intern verzoek / Actie voor medewerker"
$wsDash.Cells.Item(4, 2).Value = 1

# ---------------------------------------------------------------------
# Grow the conditional-formatting ranges on Logs from row 2:4 to 2:5 so
# the new row is covered by the same highlighting rules.
# ---------------------------------------------------------------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $fcs = $wsLogs.Range("$col" + "2:" + "$col" + "4").FormatConditions
    $newRange = $wsLogs.Range("$col" + "2:" + "$col" + "5")
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
